# Weekly update for "Fruta, Agrícola del Norte S.A. de Arica - Uva"
# - rows 79, 80, 81 get overwritten with a newer week's readings (date 44615)
# - row 82 keeps its old date (44258) but its quality/volume/price move to
#   what used to be row 79's values, while its old "second quality" figures
#   (the ones that used to live in rows 80-82 before the edit) are pushed
#   down into three brand-new rows 83, 84 and 85.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# ---- Row 79: Flame Seedless, Primera ----
$ws.Range("D79").Value = 44615
$ws.Range("K79").Value = "Flame Seedless"
$ws.Range("M79").Value = 250
$ws.Range("N79").Value = 17000
$ws.Range("O79").Value = 18000
$ws.Range("P79").Value = 17500
$ws.Range("S79").Value = 972

# ---- Row 80: Rosada pastilla, now Primera / $/caja 12 kilos ----
$ws.Range("D80").Value = 44615
$ws.Range("L80").Value = "Primera"
$ws.Range("M80").Value = 300
$ws.Range("N80").Value = 15000
$ws.Range("O80").Value = 16000
$ws.Range("P80").Value = 15500
$ws.Range("Q80").Value = "$/caja 12 kilos"
$ws.Range("S80").Value = 1292
$ws.Range("T80").Value = 12

# ---- Row 81: Superior Seedless, Primera ----
$ws.Range("D81").Value = 44615
$ws.Range("K81").Value = "Superior Seedless"
$ws.Range("L81").Value = "Primera"
$ws.Range("M81").Value = 300
$ws.Range("N81").Value = 17000
$ws.Range("O81").Value = 18000
$ws.Range("P81").Value = 17500
$ws.Range("S81").Value = 972

# ---- Row 82: Red Globe, now Primera, older date, updated volume ----
$ws.Range("D82").Value = 44258
$ws.Range("L82").Value = "Primera"
$ws.Range("M82").Value = 380

# ---- New row 83: Rosada pastilla, Segunda (previously row 80's data) ----
$ws.Range("A83").Value = 1
$ws.Range("B83").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C83").Value = "Arica y Parinacota"
$ws.Range("D83").Value = 44258
$ws.Range("D83").NumberFormat = $dateFormat
$ws.Range("E83").Value = 15
$ws.Range("F83").Value = "Fruta"
$ws.Range("G83").Value = 100109
$ws.Range("H83").Value = "Uva"
$ws.Range("I83").Value = 100109001
$ws.Range("J83").Value = "Uva"
$ws.Range("K83").Value = "Rosada pastilla"
$ws.Range("L83").Value = "Segunda"
$ws.Range("M83").Value = 370
$ws.Range("N83").Value = 19000
$ws.Range("O83").Value = 20000
$ws.Range("P83").Value = 19500
$ws.Range("Q83").Value = "$/bandeja 18 kilos"
$ws.Range("R83").Value = "Región de Coquimbo"
$ws.Range("S83").Value = 1083
$ws.Range("T83").Value = 18

# ---- New row 84: Thompson seedless, Segunda (previously row 81's data) ----
$ws.Range("A84").Value = 1
$ws.Range("B84").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C84").Value = "Arica y Parinacota"
$ws.Range("D84").Value = 44258
$ws.Range("D84").NumberFormat = $dateFormat
$ws.Range("E84").Value = 15
$ws.Range("F84").Value = "Fruta"
$ws.Range("G84").Value = 100109
$ws.Range("H84").Value = "Uva"
$ws.Range("I84").Value = 100109001
$ws.Range("J84").Value = "Uva"
$ws.Range("K84").Value = "Thompson seedless"
$ws.Range("L84").Value = "Segunda"
$ws.Range("M84").Value = 400
$ws.Range("N84").Value = 14000
$ws.Range("O84").Value = 15000
$ws.Range("P84").Value = 14500
$ws.Range("Q84").Value = "$/bandeja 18 kilos"
$ws.Range("R84").Value = "Región de Coquimbo"
$ws.Range("S84").Value = 806
$ws.Range("T84").Value = 18

# ---- New row 85: Red Globe, Segunda (previously row 82's data) ----
$ws.Range("A85").Value = 1
$ws.Range("B85").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C85").Value = "Arica y Parinacota"
$ws.Range("D85").Value = 44349
$ws.Range("D85").NumberFormat = $dateFormat
$ws.Range("E85").Value = 15
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100109
$ws.Range("H85").Value = "Uva"
$ws.Range("I85").Value = 100109001
$ws.Range("J85").Value = "Uva"
$ws.Range("K85").Value = "Red Globe"
$ws.Range("L85").Value = "Segunda"
$ws.Range("M85").Value = 300
$ws.Range("N85").Value = 11000
$ws.Range("O85").Value = 12000
$ws.Range("P85").Value = 11500
$ws.Range("Q85").Value = "$/bandeja 18 kilos"
$ws.Range("R85").Value = "Región de Coquimbo"
$ws.Range("S85").Value = 639
$ws.Range("T85").Value = 18
